$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4: "A" -> "ShadowStrike", with new stat values
$ws.Range("B4").Value = "ShadowStrike"
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 30
$ws.Range("E4").Value = 120

# Update the selection to match the author's final cursor position
$ws.Range("E11").Select()
